# Weekly crime data refresh — CompStat_1 sheet
# Updates the "Volume/Number" header, the reporting week date range,
# and the full Murder..TOTAL/Transit/Housing data table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 30   Number  28" -> "...Number  29" -------------------
$hdr = $ws.Range("A8")
$hdrText = $hdr.Value()
$oldNum = "28"
$newNum = "29"
$pos = $hdrText.IndexOf($oldNum, 10) + 1
$hdr.Characters($pos, $oldNum.Length).Text = $newNum

# --- Reporting week: 7/10/2023-7/16/2023 -> 7/17/2023-7/23/2023 ------------
$week = $ws.Range("C9")
$weekText = $week.Value()
$oldStart = "7/10/2023"
$newStart = "7/17/2023"
$oldEnd = "7/16/2023"
$newEnd = "7/23/2023"
$startPos = $weekText.IndexOf($oldStart) + 1
$week.Characters($startPos, $oldStart.Length).Text = $newStart
$weekText2 = $week.Value()
$endPos = $weekText2.IndexOf($oldEnd) + 1
$week.Characters($endPos, $oldEnd.Length).Text = $newEnd

# --- Data table (rows 14-30, columns C-N) -----------------------------------
$updates = @{
    "C14" = 6
    "D14" = 12
    "E14" = -50
    "F14" = 31
    "H14" = -35.416666666666
    "I14" = 228
    "J14" = 257
    "K14" = -11.284046692607
    "L14" = -10.9375
    "M14" = -22.972972972973
    "N14" = -78.531073446327
    "C15" = 27
    "D15" = 40
    "E15" = -32.5
    "F15" = 119
    "G15" = 136
    "H15" = -12.5
    "I15" = 825
    "J15" = 926
    "K15" = -10.907127429805
    "L15" = 0.364963503649
    "M15" = 18.364418938307
    "N15" = -55.114254624592
    "C16" = 352
    "D16" = 399
    "E16" = -11.779448621553
    "F16" = 1361
    "G16" = 1601
    "H16" = -14.990630855715
    "I16" = 8898
    "J16" = 9491
    "K16" = -6.248024444210
    "L16" = 31.355181576616
    "M16" = -12.404016538688
    "N16" = -80.887533293238
    "C17" = 573
    "D17" = 589
    "E17" = -2.716468590831
    "F17" = 2368
    "G17" = 2374
    "H17" = -0.252737994945
    "I17" = 15365
    "J17" = 14563
    "K17" = 5.507107052118
    "L17" = 27.004463547693
    "M17" = 60.873206994032
    "N17" = -33.674350341017
    "C18" = 256
    "D18" = 329
    "E18" = -22.188449848024
    "F18" = 984
    "G18" = 1222
    "H18" = -19.476268412438
    "I18" = 7697
    "J18" = 8664
    "K18" = -11.161126500461
    "L18" = 19.111730114515
    "M18" = -22.166043078167
    "N18" = -86.065499574560
    "C19" = 978
    "D19" = 1072
    "E19" = -8.768656716417
    "F19" = 3860
    "G19" = 4218
    "H19" = -8.487434803224
    "I19" = 27393
    "J19" = 28046
    "K19" = -2.328317763673
    "L19" = 45.847087637099
    "M19" = 36.195495450703
    "N19" = -40.821793514657
    "C20" = 376
    "D20" = 243
    "E20" = 54.732510288065
    "F20" = 1354
    "G20" = 1140
    "H20" = 18.771929824561
    "I20" = 8627
    "J20" = 7319
    "K20" = 17.871293892608
    "L20" = 70.831683168316
    "M20" = 52.098025387870
    "N20" = -86.057373737373
    "C21" = 2568
    "D21" = 2684
    "E21" = -4.321907600596
    "F21" = 10077
    "G21" = 10739
    "H21" = -6.164447341465
    "I21" = 69033
    "J21" = 69266
    "K21" = -0.336384373285
    "L21" = 37.395509911631
    "M21" = 22.45104299702
    "N21" = -70.751579295238
    "C22" = 28
    "D22" = 39
    "E22" = -28.205128205128
    "F22" = 162
    "G22" = 179
    "H22" = -9.497206703910
    "I22" = 1218
    "J22" = 1292
    "K22" = -5.727554179566
    "L22" = 46.043165467625
    "M22" = 5.090595340811
    "C23" = 131
    "D23" = 134
    "E23" = -2.238805970149
    "F23" = 512
    "G23" = 491
    "H23" = 4.276985743380
    "I23" = 3459
    "J23" = 3333
    "K23" = 3.780378037803
    "L23" = 17.254237288135
    "M23" = 52.177738671359
    "C24" = 2311
    "D24" = 2310
    "E24" = 0.043290043290
    "F24" = 8815
    "G24" = 9301
    "H24" = -5.225244597355
    "I24" = 61068
    "J24" = 62574
    "K24" = -2.406750407517
    "L24" = 38.999408203213
    "M24" = 38.746762393783
    "C25" = 843
    "D25" = 835
    "E25" = 0.958083832335
    "F25" = 3563
    "G25" = 3521
    "H25" = 1.192842942345
    "I25" = 24347
    "J25" = 23353
    "K25" = 4.256412452361
    "L25" = 29.691578330581
    "M25" = -6.526663339348
    "C26" = 45
    "E26" = -18.181818181818
    "F26" = 185
    "G26" = 204
    "H26" = -9.313725490196
    "I26" = 1370
    "J26" = 1504
    "K26" = -8.909574468085
    "L26" = 1.406365655070
    "C27" = 112
    "D27" = 113
    "E27" = -0.884955752212
    "F27" = 441
    "G27" = 430
    "H27" = 2.558139534883
    "I27" = 2943
    "J27" = 2857
    "K27" = 3.010150507525
    "L27" = 14.9609375
    "C28" = 27
    "D28" = 49
    "E28" = -44.897959183673
    "F28" = 127
    "G28" = 197
    "H28" = -35.532994923857
    "I28" = 671
    "J28" = 936
    "K28" = -28.311965811965
    "L28" = -32.967032967033
    "M28" = -31.390593047034
    "N28" = -79.561376789521
    "C29" = 21
    "D29" = 38
    "E29" = -44.736842105263
    "F29" = 100
    "G29" = 158
    "H29" = -36.708860759493
    "I29" = 565
    "J29" = 774
    "K29" = -27.002583979328
    "L29" = -34.225844004656
    "M29" = -29.463171036204
    "N29" = -80.950775455158
    "C30" = 3
    "D30" = 13
    "E30" = -76.923076923076
    "F30" = 28
    "G30" = 43
    "H30" = -34.883720930232
    "I30" = 268
    "J30" = 380
    "K30" = -29.473684210526
    "L30" = -15.457413249211
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
